# "divide variable beta 4" - update operand fields for instruction rows 11-13
# on Sheet1 (the binary instruction table), then move the selection to B15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 11: operand A bit-field columns D,F,G and operand B bit-field columns I,K,L
# (leading apostrophe preserves the existing quote-prefixed text style on this row)
$ws.Range("D11").Value = "'011"
$ws.Range("F11").Value = "'1"
$ws.Range("G11").Value = "'0"
$ws.Range("I11").Value = "'001"
$ws.Range("K11").Value = "'0"
$ws.Range("L11").Value = "'1"

# Row 12: operand A bit-field columns D,E
$ws.Range("D12").Value = "011"
$ws.Range("E12").Value = "1"

# Row 13: instruction column B, operand A columns D,F, operand B columns I,J,K,
# operand C columns N,O,P
$ws.Range("B13").Value = "00011"
$ws.Range("D13").Value = "011"
$ws.Range("F13").Value = "1"
$ws.Range("I13").Value = "001"
$ws.Range("J13").Value = "0"
$ws.Range("K13").Value = "1"
$ws.Range("N13").Value = "010"
$ws.Range("O13").Value = "1"
$ws.Range("P13").Value = "0"

# Move active selection to B15 (single cell), matching the saved view state.
$ws.Activate()
$ws.Range("B15").Select()
